# Nuevo formato 15 jun 2021
# Rebuild the "Rescatables" sheet with the updated remedial-exam roster:
# a few students were swapped/added and the list now runs through row 27.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$data = @(
  @(19330051920046, "COSCAHUA", "TZOYONTLE", "ALMA LIZETH", "ECOLOGÍA", "4AEV", 2),
  @(19330051920057, "GONZALEZ", "SANCHEZ", "JONATHAN", "ECOLOGÍA", "4AEV", 2),
  @(19330051920064, "LARA", "VILLAR", "VICTOR SAUL", "ECOLOGÍA", "4AEV", 2),
  @(19330051920075, "ROMERO", "CORTES", "ARTURO", "ECOLOGÍA", "4AEV", 2),
  @(19330051920418, "SAN JUAN", "CANSECO", "MARTI NEFTALI", "ECOLOGÍA", "4AEV", 2),
  @(19330051920151, "CABRERA", "GARCIA", "AYELEN", "ECOLOGÍA", "4ALCM", 2),
  @(19330051920292, "RIVERA", "FLORES", "KARLA", "ECOLOGÍA", "4ALCV", 2),
  @(19330051920295, "SANCHEZ", "TEZOCO", "ESMERALDA", "ECOLOGÍA", "4ALCV", 2),
  @(19330051920443, "TEXCAHUA", "DE LA LUZ", "ASHLEY ZURELY", "ECOLOGÍA", "4ALCV", 2),
  @(19330051920233, "GOMEZ", "MORALES", "URIEL", "ECOLOGÍA", "4APM", 2),
  @(19330051920429, "HERNANDEZ", "CALPULALPAN", "YARELY JACQUELINE", "ECOLOGÍA", "4APM", 2),
  @(19330051920235, "HERAS", "LOPEZ", "CESAR ENRIQUE", "ECOLOGÍA", "4APM", 2),
  @(19330051920139, "SOLIS", "ORTIZ", "DANIELA", "ECOLOGÍA", "4ARHM", 2),
  @(19330051920360, "BRETON", "VICENTE", "AMYRA NAHOMY", "ECOLOGÍA", "4ARHV", 2),
  @(19330051920367, "CRUZ", "LOPEZ", "ZURI SADAY", "ECOLOGÍA", "4ARHV", 2),
  @(19330051920383, "MENDEZ", "SANTOS", "FATIMA", "ECOLOGÍA", "4ARHV", 2),
  @(19330051920045, "CHAVEZ", "DE LOS SANTOS", "EUSEBIO", "ECOLOGÍA", "4AEV", 1),
  @(19330051920067, "MARTINEZ", "MONTERO", "ALEXIS YAIR", "ECOLOGÍA", "4AEV", 1),
  @(19330051920069, "MEDRANO", "LOZANO", "JOSE DANIEL", "ECOLOGÍA", "4AEV", 1),
  @(19330051920074, "REYES", "DE LA CRUZ", "IVAN", "ECOLOGÍA", "4AEV", 1),
  @(19330051920081, "VAZQUEZ", "ROMERO", "MONSERRAT", "ECOLOGÍA", "4AEV", 1),
  @(19330051920278, "GARCIA", "LINARES", "ANDRES", "ECOLOGÍA", "4ALCV", 1),
  @(19330051920436, "RAMOS", "PEREZ", "ASTRID", "ECOLOGÍA", "4ARHV", 1),
  @(19330051920369, "DE LOS SANTOS", "XOTLANIHUA", "JENNIFER", "ECOLOGÍA", "4ARHV", 1),
  @(19330051920398, "VERA", "MORALES", "VALERIA", "ECOLOGÍA", "4ARHV", 1),
  @(19330051920201, "JUAREZ", "MORO", "DENISSE", "ECOLOGÍA", "4BLCM", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}
